# Add the 2022-Q3 fund-holding sheet (new quarter) and update the summary
# ("总计") sheet with a corresponding new row, per the "feat: add 2022-Q3
# data" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计" (i.e. right
#    before the current first quarterly sheet, "2022-Q2").
# ---------------------------------------------------------------------
$summary   = $wb.Worksheets.Item("总计")
$oldFirstQ = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Add($oldFirstQ)
$q3.Name = "2022-Q3"

# Header row for the new quarterly sheet (same layout as the other
# quarterly fund-holding sheets).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q3.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 2]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows.
$q3Rows = @(
    @("210009", "金鹰核心资源混合", "2.84", "93.42", "4.95", "0.1406", 10),
    @("162102", "金鹰中小盘精选混合", "3.17", "78.28", "4.39", "0.1392", 9)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $r = $i + 2
    $row = $q3Rows[$i]

    $idx = $q3.Cells.Item($r, 1)
    $idx.Value = $i
    $idx.Font.Bold = $true
    $idx.HorizontalAlignment = -4108
    $idx.VerticalAlignment = -4160
    $idx.Borders.LineStyle = 1

    $q3.Cells.Item($r, 2).Value = "'" + $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    $q3.Cells.Item($r, 4).Value = "'" + $row[2]
    $q3.Cells.Item($r, 5).Value = "'" + $row[3]
    $q3.Cells.Item($r, 6).Value = "'" + $row[4]
    $q3.Cells.Item($r, 7).Value = "'" + $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------
# 2) Add the matching row to the "总计" summary sheet: insert a new row
#    2 with the 2022-Q3 totals, pushing the existing rows down.
# ---------------------------------------------------------------------
$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()

$a2 = $summary.Range("A2")
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.28

# Renumber the index column (A) for the rows that got pushed down so it
# stays a contiguous 0..5 sequence.
for ($r = 3; $r -le 7; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# Restore the originally active sheet ("2020-Q4" was the last/active tab)
# so adding the new sheet doesn't otherwise disturb workbook UI state.
$wb.Worksheets.Item("2020-Q4").Activate()
